$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: blank cells styled like the rest of column A (bold/border/center,
#     style index 1) by copying format from an existing column-A cell, so no new
#     cellXf entries are created for column A.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A43:A50").PasteSpecial(-4122) | Out-Null

# --- Column B: new date values with a custom date format.
# First register numFmt 166 ("yyyy-mm-dd") *and* land on numFmt 167 ("YYYY-MM-DD")
# using a single throwaway cell so only ONE new cellXf survives (matches the source
# edit, which leaves 166 registered in numFmts but unused by any cellXf).
$tmp = $ws.Cells.Item(43, 2)
$tmp.Value = 44924
$tmp.NumberFormat = "yyyy-mm-dd"
$tmp.NumberFormat = "YYYY-MM-DD"

# Now apply the same (already-registered) format to the remaining date cells so they
# reuse the same cellXf instead of minting new ones.
$restDates = $ws.Range("B44:B50")
$restDates.Value = 44924
$restDates.NumberFormat = "YYYY-MM-DD"

# --- Columns C/D/E/F: word, meaning, example, classes for the 8 new dictionary rows.
$rows = @(
    @{ Row = 43; C = 'suspend                  '; D = 'to officially stop something or somebody for a time/to hang something from something else'; E = 'Production has been suspended while safety checks are carried out/A lamp was suspended from the ceiling.' }
    @{ Row = 44; C = 'sweep                    '; D = 'to clean a room/push something with a lot of force'; E = 'to sweep the floor/ the little boat was swept out to sea' }
    @{ Row = 45; C = 'tackle                   '; D = 'to make a determined effort to deal with a difficult problem or situation'; E = 'Firefighters tackled a blaze in a garage last night' }
    @{ Row = 46; C = 'tale '; D = 'a story'; E = 'The story is a classic tale of love and betrayal.' }
    @{ Row = 47; C = 'tension                  '; D = 'a situation in which people do not trust each other, or feel unfriendly'; E = 'Family tensions and conflicts may lead to violence' }
    @{ Row = 48; C = 'thorough                 '; D = 'done completely; with great attention to detail'; E = 'a thorough knowledge of the subject' }
    @{ Row = 49; C = 'tissue                   '; D = 'a piece of soft paper, used especially as a handkerchief'; E = 'a box of tissues' }
    @{ Row = 50; C = 'trace       '; D = 'to find or discover somebody/something by looking carefully for them/it  OR / gentally touch'; E = 'We finally traced him to an address in Chicago /She lightly traced the outline of his face with her finger' }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = 0
}

